$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the song notes text for lostTogether (row 3, column C):
# "quint 91 kaleids" -> "2 pairs split 2/3 on M1/M2"
$ws.Range("C3").Value = "~cue 9: 2 x double heavens gate                                                                                            ~cues 10, 11: 2 pairs split 2/3 on M1/M2"

# Update numCues for lostTogether (row 3, column B) from 40 to 42
$ws.Range("B3").Value = 42

# Recalculate so the SUM formula in B8 updates
$excel.Calculate()

# Update the selected cell in the sheet view to C3
$ws.Range("C3").Select()
